$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format before writing, to prevent Excel from
# auto-converting numeric-looking strings (e.g. "240.80") into real numbers.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range('D2').Value = '43.898.10'
$ws.Range('E2').Value = '  -0.99%  '
$ws.Range('D3').Value = '2.350.95'
$ws.Range('E3').Value = '  -0.87%  '
$ws.Range('E5').Value = '  -2.86%  '
$ws.Range('D6').Value = '240.80'
$ws.Range('E6').Value = '  -1.62%  '
$ws.Range('D7').Value = '72.89'
$ws.Range('E7').Value = '  -5.01%  '
$ws.Range('D9').Value = '0.602'
$ws.Range('E9').Value = '  +1.08%  '
$ws.Range('E10').Value = '  -2.45%  '
$ws.Range('D11').Value = '59.24'
$ws.Range('E11').Value = '  +2.09%  '
$ws.Range('D12').Value = '33.26'
$ws.Range('E12').Value = '  +2.35%  '
$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D13').Value = '0.109'
$ws.Range('E13').Value = '  -0.17%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').Value = '7.30'
$ws.Range('E14').Value = '  -2.63%  '
$ws.Range('D15').Value = '2.704.08'
$ws.Range('E15').Value = '  -0.78%  '
$ws.Range('D16').Value = '16.38'
$ws.Range('E16').Value = '  -4.65%  '
$ws.Range('D17').Value = '0.907'
$ws.Range('E17').Value = '  -1.92%  '
$ws.Range('D18').Value = '2.352.41'
$ws.Range('E18').Value = '  -1.20%  '
$ws.Range('D19').Value = '43.814.58'
$ws.Range('E19').Value = '  -1.52%  '
$ws.Range('D20').Value = '0.0000104'
$ws.Range('E20').Value = '  -0.21%  '
$ws.Range('E21').Value = '  -0.13%  '
$ws.Range('D22').Value = '77.78'
$ws.Range('E22').Value = '  -1.09%  '
$ws.Range('D23').Value = '256.08'
$ws.Range('E23').Value = '  -1.12%  '
$ws.Range('D24').Value = '1.95'
$ws.Range('E24').Value = '  +12.99%  '
$ws.Range('E25').Value = '  -0.01%  '
$ws.Range('D26').Value = '3.74'
$ws.Range('E26').Value = '  +0.89%  '
$ws.Range('D27').Value = '2.50'
$ws.Range('E27').Value = '  -3.37%  '
$ws.Range('D28').Value = '10.56'
$ws.Range('E28').Value = '  -2.89%  '
$ws.Range('D29').Value = '2.28'
$ws.Range('E29').Value = '  -1.72%  '
$ws.Range('D30').Value = '22.58'
$ws.Range('E30').Value = '  -2.15%  '
$ws.Range('E31').Value = '  +1.04%  '
$ws.Range('E32').Value = '  -1.41%  '
$ws.Range('E33').Value = '  +0.75%  '
$ws.Range('D34').Value = '0.0753'
$ws.Range('E34').Value = '  -1.44%  '
$ws.Range('B35').Value = 'InternetComputer(DFINITY)'
$ws.Range('C35').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D35').Value = '5.46'
$ws.Range('E35').Value = '  +1.80%  '
$ws.Range('B36').Value = 'Filecoin'
$ws.Range('C36').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D36').Value = '5.12'
$ws.Range('E36').Value = '  -5.77%  '
$ws.Range('E37').Value = '  -2.85%  '
$ws.Range('D38').Value = '6.43'
$ws.Range('E38').Value = '  -2.82%  '
$ws.Range('D39').Value = '2.38'
$ws.Range('E39').Value = '  -4.63%  '
$ws.Range('E40').Value = '  -0.16%  '
$ws.Range('E41').Value = '  +25.75%  '
$ws.Range('E42').Value = '  +14.72%  '
$ws.Range('D43').Value = '0.109'
$ws.Range('E43').Value = '  +7.76%  '
$ws.Range('D44').Value = '9.27'
$ws.Range('E44').Value = '  +1.16%  '
$ws.Range('B45').Value = 'InjectiveProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D45').Value = '18.99'
$ws.Range('E45').Value = '  -1.39%  '
$ws.Range('B46').Value = 'Algorand'
$ws.Range('C46').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D46').Value = '0.202'
$ws.Range('E46').Value = '  +3.49%  '
$ws.Range('E47').Value = '  -0.93%  '
$ws.Range('D48').Value = '1.25'
$ws.Range('E48').Value = '  -1.77%  '
$ws.Range('B50').Value = 'Aave'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D50').Value = '99.59'
$ws.Range('E50').Value = '  -2.52%  '
$ws.Range('B51').Value = 'ARBITRUM'
$ws.Range('C51').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D51').Value = '1.16'
$ws.Range('E51').Value = '  -4.86%  '

# Restore the original (default/General) formatting so the saved file
# has no residual style deltas versus the source workbook.
$dRange.NumberFormat = "General"
$dRange.Style = "Normal"
